function Set-RowRange($WS, $Row, $StartCol, $Values) {
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $WS.Cells.Item($Row, $StartCol + $i).Value = $Values[$i]
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scalar parameter changes (rows 2, 3, 6) ---
$ws.Range("C2").Value = "Contrachapado"
$ws.Range("C3").Value = 5
$ws.Range("C6").Value = 174

# --- Data table changes (rows 9-12, columns C:AG) ---
$row9 = @(14.58, 16.52, 18.53, 20.6, 22.54, 24.55, 26.62, 28.56, 30.5, 32.65, 15.07, 21.59, 26.22, 30.29, 33.76, 37.18, 40.59, 43.69, 46.75, 50.08, 53.07, 56.03, 59.09, 62.23, 65.16, 68.19, 71.31, 74.23, 77.14, 80.36, 83.26)
Set-RowRange $ws 9 3 $row9
$row10 = @(21.5, 22, 22.9, 24.1, 25.2, 26.5, 27.9, 29.1, 30, 27.9, 22.8, 24.5, 28.7, 32.5, 35.8, 39.1, 42.4, 45.4, 48.3, 51.5, 54.3, 57.1, 59.9, 62.7, 65.3, 67.8, 70.3, 72.5, 74.6, 76.7, 78.6)
Set-RowRange $ws 10 3 $row10
$row11 = @(13.93, 15.85, 17.85, 19.92, 21.85, 23.85, 25.93, 25.89, 25.85, 25.79, 25.72, 28.56, 31.5, 34.55, 37.41, 40.37, 43.44, 46.31, 49.19, 52.37, 55.25, 58.13, 61.12, 64.21, 67.1, 70.1, 73.2, 76.09, 78.99, 82.19, 85.09)
Set-RowRange $ws 11 3 $row11
$row12 = @(17.74, 19.02, 20.35, 21.7, 22.91, 24.02, 24.83, 24.81, 22.78, 15.91, 19.7, 23.38, 28.01, 32.07, 35.55, 38.97, 42.38, 45.48, 48.54, 51.87, 54.85, 57.82, 60.88, 64.02, 66.95, 69.98, 73.1, 76.02, 78.93, 82.15, 85.05)
Set-RowRange $ws 12 3 $row12

Write-Output "edit applied"
